# Tab25 update: South Sudan becomes resource-rich ("South Sudan*"),
# Nigeria loses its resource-rich marker ("Nigeria*" -> "Nigeria"),
# row shading follows the asterisk flag, and the dependent aggregate/
# regional-grouping rows are refreshed with their recalculated totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab25")

# --- Country-name / resource-rich marker updates -------------------------
$ws.Range("B34").Value = "South Sudan*"
$ws.Range("B57").Value = "Nigeria"

# --- Row shading follows the resource-rich marker -------------------------
# South Sudan (row 34) is now resource-rich: copy the shaded format used by
# another resource-rich country row (Chad*, row 17).
$ws.Range("B17:M17").Copy() | Out-Null
$ws.Range("B34:M34").PasteSpecial(-4122) | Out-Null

# Nigeria (row 57) is no longer resource-rich: copy the plain (unshaded)
# format used by another non-resource-rich country row (Eswatini, row 5).
$ws.Range("B5:M5").Copy() | Out-Null
$ws.Range("B57:M57").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Recalculated aggregate / regional-grouping rows ----------------------
$ws.Range("C69").Value = 59157.599999999999
$ws.Range("E69").Value = 6732.6
$ws.Range("F69").Value = 25377.1
$ws.Range("G69").Value = 20259.3
$ws.Range("I69").Value = 1219.0999999999999
$ws.Range("L69").Value = 24410.5
$ws.Range("M69").Value = 4239.8999999999996
$ws.Range("C77").Value = 830120.7
$ws.Range("D77").Value = 56200.9
$ws.Range("E77").Value = 147189.4
$ws.Range("F77").Value = 340511.9
$ws.Range("G77").Value = 64560.7
$ws.Range("H77").Value = 9861.9
$ws.Range("L77").Value = 317014.7
$ws.Range("M77").Value = 171566.2
$ws.Range("C80").Value = 211076.5
$ws.Range("D80").Value = 13983.8
$ws.Range("E80").Value = 20366.400000000001
$ws.Range("F80").Value = 70169
$ws.Range("G80").Value = 32152.799999999999
$ws.Range("H80").Value = 2285.6999999999998
$ws.Range("I80").Value = 446.7
$ws.Range("K80").Value = 4638.8999999999996
$ws.Range("L80").Value = 67581.600000000006
$ws.Range("M80").Value = 67032.5
$ws.Range("C82").Value = 1039600.2
$ws.Range("D82").Value = 83573.8
$ws.Range("E82").Value = 127371.6
$ws.Range("F82").Value = 272658.7
$ws.Range("G82").Value = 135934.1
$ws.Range("H82").Value = 11912.3
$ws.Range("I82").Value = 14643.6
$ws.Range("K82").Value = 5771.4
$ws.Range("L82").Value = 262639.7
$ws.Range("M82").Value = 387193.7
$ws.Range("C84").Value = 116847.6
$ws.Range("D84").Value = 672.9
$ws.Range("E84").Value = 15550.9
$ws.Range("F84").Value = 45945.599999999999
$ws.Range("G84").Value = 32848.9
$ws.Range("H84").Value = 2409.3000000000002
$ws.Range("I84").Value = 1584.3
$ws.Range("K84").Value = 1866.5
$ws.Range("L84").Value = 44570.5
$ws.Range("M84").Value = 15857.9
$ws.Range("C86").Value = 679975.7
$ws.Range("D86").Value = 39993.699999999997
$ws.Range("E86").Value = 77890
$ws.Range("F86").Value = 227539.1
$ws.Range("G86").Value = 119484.8
$ws.Range("H86").Value = 6315.9
$ws.Range("I86").Value = 8437.7000000000007
$ws.Range("K86").Value = 5856
$ws.Range("L86").Value = 220356
$ws.Range("M86").Value = 194217.2
$ws.Range("C87").Value = 4430047.9000000004
$ws.Range("E87").Value = 981599.9
$ws.Range("F87").Value = 767127.1
$ws.Range("G87").Value = 387392.2
$ws.Range("H87").Value = 78037.899999999994
$ws.Range("K87").Value = 45546.2
$ws.Range("L87").Value = 710764.5
$ws.Range("M87").Value = 1967633.5
$ws.Range("C89").Value = 14397124.1
$ws.Range("D89").Value = 536319.80000000005
$ws.Range("E89").Value = 3593420.4
$ws.Range("F89").Value = 1936951.7
$ws.Range("G89").Value = 756839.2
$ws.Range("H89").Value = 203681.1
$ws.Range("I89").Value = 181899.5
$ws.Range("J89").Value = 4494.5
$ws.Range("K89").Value = 78327.100000000006
$ws.Range("L89").Value = 1634659.6
$ws.Range("M89").Value = 7105202
$ws.Range("C90").Value = 10761783.699999999
$ws.Range("D90").Value = 744101.4
$ws.Range("E90").Value = 1418269.5
$ws.Range("F90").Value = 3212470.4
$ws.Range("G90").Value = 847869.2
$ws.Range("H90").Value = 481725
$ws.Range("I90").Value = 136249.29999999999
$ws.Range("J90").Value = 13322.8
$ws.Range("K90").Value = 19048.599999999999
$ws.Range("L90").Value = 2879945.9
$ws.Range("M90").Value = 3888716.2
$ws.Range("C94").Value = 134730.29999999999
$ws.Range("D94").Value = 10607.1
$ws.Range("E94").Value = 26677.599999999999
$ws.Range("F94").Value = 28042.3
$ws.Range("G94").Value = 4901.7
$ws.Range("H94").Value = 1255.0999999999999
$ws.Range("L94").Value = 24485.1
$ws.Range("M94").Value = 59392.2
$ws.Range("C97").Value = 392794.4
$ws.Range("E97").Value = 39040.1
$ws.Range("F97").Value = 158338.20000000001
$ws.Range("G97").Value = 97349.6
$ws.Range("H97").Value = 6502.7
$ws.Range("I97").Value = 3286
$ws.Range("L97").Value = 154993.1
$ws.Range("M97").Value = 69186.600000000006
$ws.Range("C98").Value = 1286621.1000000001
$ws.Range("D98").Value = 78152.3
$ws.Range("E98").Value = 235885.6
$ws.Range("F98").Value = 306949
$ws.Range("G98").Value = 185299.4
$ws.Range("H98").Value = 40111.599999999999
$ws.Range("K98").Value = 19670.7
$ws.Range("L98").Value = 286562.90000000002
$ws.Range("M98").Value = 399272.9
